# Update metrics values for rows 2-26 (columns B:I) with new rstd-based results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(0.2937639237727891, 0.3846156044978735, -1.099699748824124, 0.1626454243347781, 0.7815952301025391, 1.112792730331421, 1.066634297370911, 1.091071844100952)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 9; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
